$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style used by the existing header row (s="1", e.g. AC1)
# onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every team in this sheet shares the same 2022 Oakland A's season record:
# 60 wins, 102 losses, 0 ties. Fill it in for every data row (2-66).
for ($r = 2; $r -le 66; $r++) {
    $ws.Cells.Item($r, 30).Value = 60
    $ws.Cells.Item($r, 31).Value = 102
    $ws.Cells.Item($r, 32).Value = 0
}
